# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price cells that look like plain numbers ("6.09", "154.28", ...) are forced
# back to text via a temporary "@" (Text) number format so Excel doesn't
# silently coerce them to numeric values; the style is reset to "Normal"
# right after so no stray cell-style index is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.302.83"
$ws.Range("E2").Value = "  +2.39%  "

$ws.Range("D3").Value = "3.203.32"
$ws.Range("E3").Value = "  +1.89%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("E5").Value = "  +1.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.79%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "3.202.42"
$ws.Range("E8").Value = "  +1.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.40%  "

$ws.Range("E10").Value = "  +4.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.75%  "

$ws.Range("E13").Value = "  +3.20%  "

$ws.Range("E14").Value = "  +6.23%  "

$ws.Range("D15").Value = "3.734.87"
$ws.Range("E15").Value = "  +1.92%  "

$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("E17").Value = "  +4.75%  "

$ws.Range("D18").Value = "65.047.24"
$ws.Range("E18").Value = "  +2.28%  "

$ws.Range("D19").Value = "3.203.00"
$ws.Range("E19").Value = "  +1.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.775"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "83.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.87%  "

$ws.Range("E27").Value = "  +0.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.72%  "

$ws.Range("E29").Value = "  +3.72%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.20%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.89%  "

$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("E33").Value = "  +9.22%  "

$ws.Range("E34").Value = "  +5.73%  "

$ws.Range("D35").Value = "0.0₃0904"
$ws.Range("E35").Value = "  +5.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.80%  "

$ws.Range("E37").Value = "  +4.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.49%  "

$ws.Range("E39").Value = "  +3.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "478.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "51.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.19%  "

$ws.Range("E43").Value = "  +9.09%  "

$ws.Range("E44").Value = "  +3.34%  "

$ws.Range("D45").Value = "2.962.12"
$ws.Range("E45").Value = "  +1.40%  "

$ws.Range("E46").Value = "  +3.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.53%  "

$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.79%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.35%  "

$ws.Range("E51").Value = "  +0.00%  "
